# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
# - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
# - Latest handback timestamps for zh-cn and de-de are refreshed
# - Error details are cleared now that files are back in sync

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: Status columns for zh-cn (E2) and de-de (F2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn detail sheet
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-01 12:55:43"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# de-de detail sheet
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-01 12:55:51"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
